$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.493.68"
$ws.Range("E2").Value = "  +1.19%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.865.94"
$ws.Range("E3").Value = "  +1.37%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'325.25"
$ws.Range("E5").Value = "  -0.19%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.04%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4559"
$ws.Range("E7").Value = "  -1.81%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.59%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07828"
$ws.Range("E9").Value = "  -0.35%  "

# Row 10 - Polygon
$ws.Range("D10").Value = "'0.9896"
$ws.Range("E10").Value = "  +2.75%  "

# Row 11 - Solana
$ws.Range("D11").Value = "'21.57"
$ws.Range("E11").Value = "  -2.24%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.862.83"
$ws.Range("E12").Value = "  +2.08%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "'6.911"
$ws.Range("E13").Value = "  +0.55%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.643"
$ws.Range("E14").Value = "  -0.90%  "

# Row 15 - TRON
$ws.Range("D15").Value = "'0.06966"
$ws.Range("E15").Value = "  +0.77%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'86.74"
$ws.Range("E16").Value = "  -2.01%  "

# Row 17 - BinanceUSD
$ws.Range("E17").Value = "  -0.04%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "'0.000009946"
$ws.Range("E18").Value = "  -0.10%  "

# Row 19 - Avalanche
$ws.Range("D19").Value = "'16.67"
$ws.Range("E19").Value = "  -0.28%  "

# Row 20 - Dai
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  -0.05%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "28.482.80"

# Row 22 - Uniswap
$ws.Range("D22").Value = "'5.255"
$ws.Range("E22").Value = "  -0.90%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'10.92"
$ws.Range("E23").Value = "  -1.00%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.20%  "

# Row 25 - WrappedliquidstakedEther2.0
$ws.Range("D25").Value = "2.084.13"
$ws.Range("E25").Value = "  +1.61%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'153.85"
$ws.Range("E26").Value = "  -0.37%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'19.13"
$ws.Range("E27").Value = "  -0.20%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "'5.674"
$ws.Range("E28").Value = "  -1.42%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "'1.939"
$ws.Range("E29").Value = "  -1.53%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "'117.32"
$ws.Range("E30").Value = "  -1.24%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.09272"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'0.9090"
$ws.Range("E32").Value = "  -2.31%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "'5.264"
$ws.Range("E33").Value = "  -0.49%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  -0.47%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'3.296"
$ws.Range("E35").Value = "  -1.12%  "

# Row 36 - Hedera
$ws.Range("D36").Value = "'0.05717"
$ws.Range("E36").Value = "  -1.58%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "'1.140"
$ws.Range("E37").Value = "  -0.79%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -2.06%  "

# Row 39 - FraxShare
$ws.Range("D39").Value = "'7.683"
$ws.Range("E39").Value = "  -1.09%  "

# Row 40 - TheSandbox
$ws.Range("D40").Value = "'0.5569"
$ws.Range("E40").Value = "  -0.35%  "

# Row 41 - Algorand
$ws.Range("E41").Value = "  +0.53%  "

# Row 42 - Aptos
$ws.Range("D42").Value = "'9.640"
$ws.Range("E42").Value = "  -2.44%  "

# Row 43 - Cronos
$ws.Range("D43").Value = "'0.07097"
$ws.Range("E43").Value = "  -1.85%  "

# Row 44 - EnergySwap
$ws.Range("D44").Value = "'11.60"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "'0.5239"
$ws.Range("E45").Value = "  -0.62%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "'2.142"
$ws.Range("E46").Value = "  +0.74%  "

# Row 47 - WEMIXToken
$ws.Range("D47").Value = "'1.128"
$ws.Range("E47").Value = "  -0.52%  "

# Row 48 - NEARProtocol
$ws.Range("E48").Value = "  -1.09%  "

# Row 49 - Quant
$ws.Range("D49").Value = "'111.72"
$ws.Range("E49").Value = "  -1.83%  "

# Row 50 - MXToken
$ws.Range("D50").Value = "'2.416"
$ws.Range("E50").Value = "  +4.07%  "

# Row 51 - PaxDollar
$ws.Range("D51").Value = "'1.006"
$ws.Range("E51").Value = "  -0.10%  "
